$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content corrections ---
# Row 10 (item 9, SMA/coax connector): the "References" designator changes
# from "J1" to "+"
$ws.Range("D10").Value = "+"

# Row 11 (item 10, 3-pin connector): "Description" simplified to the part
# name "Conn_01x03" (was the long KiCad-generated description)
$ws.Range("B11").Value = "Conn_01x03"

# --- New "notified" marker column (I) ---
# A new column I is populated with 1 for the rows whose notification has
# been completed.
$notifiedRows = 2,4,5,7,8,9,10,14,15,18,19,20,21,22,23,24,25
foreach ($r in $notifiedRows) {
    $ws.Cells.Item($r, 9).Value = 1
}
